# Generate Report for Handoff
# Marks the "Ready for handoff" files as handed off ("ht") and refreshes
# the handoff generation timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 12, 13)

$zhCnHandoffTime = "2016-08-22 14:22:04"
$deDeHandoffTime = "2016-08-22 14:22:18"

foreach ($r in $rows) {
    # Priority column on the language sheets: mark as handed off.
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # Latest Handoff Datetime on the language sheets.
    $wsZhCn.Range("H$r").Value = $zhCnHandoffTime
    $wsDeDe.Range("H$r").Value = $deDeHandoffTime

    # Latest HO Xliff Generate Date on the Overview sheet (matches de-de time).
    $wsOverview.Range("G$r").Value = $deDeHandoffTime
}
